# Add input file pattern support (e.g. Folder/*.xls*)
#
# Settings sheet:
#  - AppFilePath value: backslash path -> forward-slash path
#  - ErrorsFolder value: "Errors" -> "../TestData/Errors"
#  - TempFolder value: "Temp" -> "../TestData/Temp"
#  - InputFile -> InputPath, value becomes a wildcard pattern
#  - New OutputPath row inserted right after InputPath (archive folder for processed files)
#  - TransactionQueue row pushed down (now after a blank row)
#
# Dispatch sheet:
#  - OutputPath value: ".../Output/[File]" -> ".../Temp/[File]"
#  - CompleteFolder value: ".../Processed" -> ".../Output"

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$dispatch = $wb.Worksheets.Item("Dispatch")

# --- Settings sheet -------------------------------------------------

# AppFilePath: normalize path separators
$settings.Range("B9").Value = "C:/WINDOWS/system32/notepad.exe"

# ErrorsFolder / TempFolder: point at the TestData tree
$settings.Range("B12").Value = "../TestData/Errors"
$settings.Range("B13").Value = "../TestData/Temp"

# InputFile -> InputPath (wildcard pattern instead of a single sample file)
$settings.Range("A15").Value = "InputPath"
$settings.Range("B15").Value = "../TestData/*.xls*"

# Row 16 used to hold TransactionQueue; it now holds the new OutputPath
# entry. TransactionQueue moves down to row 18 (leaving row 17 blank as
# a separator, matching the spacing used between every other group on
# this sheet). Cells are written directly in place rather than via a
# sheet-wide row insert, so the unrelated rows 21/22 further down stay
# put.
$settings.Range("A16").Value = "OutputPath"
$settings.Range("B16").Value = "../TestData/Output"
$settings.Range("C16").Value = "Archive folder path for processed files"

$settings.Range("A18").Value = "TransactionQueue"
$settings.Range("B18").Value = "RFW-ToProcess"
$settings.Range("C18").Value = "Transactions queue in Orchestrator"

# --- Dispatch sheet --------------------------------------------------

$dispatch.Range("B13").Value = "../TestData/Temp/[File]"
$dispatch.Range("B16").Value = "../TestData/Output"

# --- Active sheet / selection ----------------------------------------
# Restore the Dispatch sheet's old selection, then make Settings the
# active tab with its new selection.

$dispatch.Activate() | Out-Null
$dispatch.Range("B17").Select() | Out-Null

$settings.Activate() | Out-Null
$settings.Range("B15").Select() | Out-Null
